$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.677.35"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "2.273.02"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.93%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.620"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.597"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0895"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.968"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").Value = "2.615.72"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").Value = "2.290.77"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").Value = "42.273.81"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.70%  "
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("E21").Value = "  +2.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("E23").Value = "  -6.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.66%  "
$ws.Range("E25").Value = "  -3.16%  "
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.00%  "
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +12.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "163.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("E33").Value = "  -3.12%  "
$ws.Range("E34").Value = "  -2.67%  "
$ws.Range("E35").Value = "  +1.88%  "
$ws.Range("E36").Value = "  -3.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.76%  "
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("E40").Value = "  -1.70%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "67.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("D46").Value = "1.707.34"
$ws.Range("E46").Value = "  +6.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "108.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "75.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.62%  "
$ws.Range("E51").Value = "  -2.56%  "
